# Commit: "Added multiThreading and Multiprocessing"
#
# 1) Sheet1: insert a new blank row at the top (shifts all existing data
#    down by one row: old A1 -> A2, old A229 -> A230, etc.)
# 2) Add a new "Sheet2" (placed right after Sheet1) containing the parsed
#    "avg / total ..." summary rows (10 rows x 7 columns).
# 3) Update view/selection state to match: Sheet1 scrolled down with a
#    multi-cell selection, Sheet2 active/selected with F1:F10 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) Shift Sheet1 data down by inserting a row at the top ---
$ws1.Rows.Item(1).Insert()

# --- 2) Create Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

$data = @(
    @(0.65700000000000003, 0.65900000000000003, 0.65200000000000002, 1600),
    @(0.65900000000000003, 0.66100000000000003, 0.65300000000000002, 1600),
    @(0.65700000000000003, 0.66,                 0.65200000000000002, 1600),
    @(0.66300000000000003, 0.66600000000000004, 0.65800000000000003, 1600),
    @(0.65800000000000003, 0.66100000000000003, 0.65300000000000002, 1600),
    @(0.65900000000000003, 0.66200000000000003, 0.65400000000000003, 1600),
    @(0.66200000000000003, 0.66500000000000004, 0.65800000000000003, 1600),
    @(0.66,                0.66400000000000003, 0.65600000000000003, 1600),
    @(0.66400000000000003, 0.66700000000000004, 0.65900000000000003, 1600),
    @(0.66200000000000003, 0.66400000000000003, 0.65700000000000003, 1600)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws2.Cells.Item($r, 1).Value = "avg"
    $ws2.Cells.Item($r, 2).Value = "/"
    $ws2.Cells.Item($r, 3).Value = "total"
    $ws2.Cells.Item($r, 4).Value = $data[$i][0]
    $ws2.Cells.Item($r, 5).Value = $data[$i][1]
    $ws2.Cells.Item($r, 6).Value = $data[$i][2]
    $ws2.Cells.Item($r, 7).Value = $data[$i][3]
}

$ws2.Columns.Item(1).ColumnWidth = 38

# --- 3) View / selection state ---

# Sheet1: scrolled down, multi-area selection across the repeated
# "avg / total" rows (now at rows 17,40,63,...,224 after the insert).
# A224 is listed first so it becomes the ActiveCell (matches the target
# selection's activeCell="A224" / activeCellId="9" - last area clicked).
$ws1.Activate()
$sel1 = $ws1.Range("A224,A17,A40,A63,A86,A109,A132,A155,A178,A201")
$sel1.Select()
$excel.ActiveWindow.ScrollRow = 208

# Sheet2: active tab, F1:F10 selected.
$ws2.Activate()
$ws2.Range("F1:F10").Select()

try { $excel.WindowState = -4140 } catch {}
try { $excel.ActiveWindow.WindowState = -4140 } catch {}
